$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (N) mirroring the existing 2019 column (M),
# copying formats first so the new cells match the surrounding table.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 2020

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 15

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 1308.3

# Move the active selection to the new cell, matching the saved view state.
$null = $ws.Range("N6").Select()
